# Team07Report.xlsx update - "Add files via upload"
#
# Content edits made by the author (the cosmetic row-height / column-width /
# dyDescent / window-position churn visible in the raw OOXML diff is just
# Excel's own re-layout noise from opening the file on a different machine
# and is not something this script needs to reproduce):
#
#   Burndown sheet:
#     C3  29  -> 28    (Sprint 1 "Stories Left")
#     D3        -> formula =C2-C3   (Stories completed this sprint)
#     E3  116 -> 229   (LOC)
#     F3  115 -> 190   (Minutes)
#     G3  stays =E3/F3, recalculates automatically from the new E3/F3
#
#   Sprint1 sheet:
#     Row 20 ("Marriage After 14"): mark Done, fill in actual size/time
#       D20 = "Done"
#       G20 = 35
#       H20 = 20
#       I20 = "yes"
#     Row 21 ("Corresponding Entries"): mark Done, fill in actual size/time
#       D21 = "Done"
#       G21 = 78
#       H21 = 55
#       I21 = "yes"

$wb = $excel.ActiveWorkbook

# --- Burndown sheet ---------------------------------------------------
$burndown = $wb.Worksheets.Item("Burndown")

$burndown.Range("C3").Value = 28
$burndown.Range("D3").Formula = "=C2-C3"
$burndown.Range("E3").Value = 229
$burndown.Range("F3").Value = 190

# --- Sprint1 sheet -----------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")

$sprint1.Range("D20").Value = "Done"
$sprint1.Range("G20").Value = 35
$sprint1.Range("H20").Value = 20
$sprint1.Range("I20").Value = "yes"

$sprint1.Range("D21").Value = "Done"
$sprint1.Range("G21").Value = 78
$sprint1.Range("H21").Value = 55
$sprint1.Range("I21").Value = "yes"
